# Updates the cryptocurrency Price (column D) and Volume(1h) (column E) values
# on the active worksheet to reflect the latest scrape, per the commit:
# "Updated cryptos list on Thu May 23 20:54:23 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.651.50"
$ws.Range("E2").Value = "  -2.73%  "

$ws.Range("D3").Value = "3.737.56"
$ws.Range("E3").Value = "  -0.39%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.57"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.66"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.61%  "

$ws.Range("D7").Value = "3.746.48"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.517"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.91%  "

$ws.Range("E10").Value = "  -5.61%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.15"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -6.06%  "

$ws.Range("E12").Value = "  -4.76%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.45"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -6.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000241"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -4.93%  "

$ws.Range("D15").Value = "4.359.03"
$ws.Range("E15").Value = "  -0.36%  "

$ws.Range("D16").Value = "3.728.94"
$ws.Range("E16").Value = "  -0.42%  "

$ws.Range("D17").Value = "67.620.51"
$ws.Range("E17").Value = "  -2.84%  "

$ws.Range("E18").Value = "  -5.00%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.13"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -4.04%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.11"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.54%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "487.74"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.49%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.98"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.12%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.714"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.76%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.37"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.87%  "

$ws.Range("E25").Value = "  -9.82%  "

$ws.Range("E26").Value = "  +5.65%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.12"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -5.89%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.20"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -8.21%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.07%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.94"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.22%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.38"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.84%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.28"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +6.26%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.63"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -5.41%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.108"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.66%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -3.97%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.135"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.43%  "

$ws.Range("E38").Value = "  -6.76%  "

$ws.Range("E39").Value = "  -7.26%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "449.86"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.36%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "48.84"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.73%  "

$ws.Range("E42").Value = "  -3.78%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.83"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -7.43%  "

$ws.Range("E44").Value = "  -3.45%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.17"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -7.48%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "140.81"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.65%  "

$ws.Range("D47").Value = "2.785.34"
$ws.Range("E47").Value = "  -5.56%  "

$ws.Range("E48").Value = "  +0.02%  "

$ws.Range("E49").Value = "  -3.33%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.65"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -5.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.96"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +7.75%  "
